# Sample Project rules workbook - row 11 "Rule" name changed from "R40" to "1".
# (Rule names in column B are stored as text in the sheet's shared-string
# table, so the new value must remain text even though it looks numeric.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Leading apostrophe forces Excel to keep the value as literal text ("1")
# instead of auto-converting it to the number 1.
$ws.Range("B11").Value = "'1"
